$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Fitness") holds 7293 for every data row (2..252); update to 7573.
$ws.Range("C2:C252").Value = 7573
